# Updates the cryptos price list (D = Price, E = Volume(1h)) to match the
# latest scrape. Helper that writes a value into a cell while guaranteeing
# it lands as literal text (Excel otherwise "helpfully" re-types numeric
# looking strings like "1.003" as a Double), then restores the cell's
# original (default/"Normal") style so no visible formatting changes.
function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price, Volume(1h)) for rows whose coin identity doesn't change.
$updates = [ordered]@{
    2  = @("26.010.92", "  -0.69%  ")
    3  = @("1.652.29",  "  -0.12%  ")
    4  = @("1.003",     "  -0.24%  ")
    5  = @("217.37",    "  -0.34%  ")
    6  = @("0.5204",    "  +0.07%  ")
    7  = @("1.003",     "  -0.19%  ")
    8  = @("0.2614",    "  -1.20%  ")
    9  = @("0.06268",   "  -0.67%  ")
    10 = @("20.51",     "  -3.43%  ")
    11 = @("0.07737",   "  +0.03%  ")
    12 = @("4.462",     "  +0.86%  ")
    13 = @("1.652.84",  "  -0.11%  ")
    14 = @("1.881.12",  "  +0.05%  ")
    15 = @("0.5422",    "  -0.66%  ")
    16 = @("0.0$([char]0x2085)8090", "  -1.45%  ")
    17 = @("65.00",     "  +0.44%  ")
    18 = @("26.029.53", "  -0.66%  ")
    19 = @($null,       "  -0.24%  ")
    20 = @("4.569",     "  -2.33%  ")
    21 = @("191.17",    "  -0.13%  ")
    22 = @("10.01",     "  -1.44%  ")
    23 = @("5.978",     "  -3.29%  ")
    24 = @("1.004",     "  -0.32%  ")
    25 = @("138.55",    "  +0.07%  ")
    26 = @("0.1235",    "  -0.54%  ")
    27 = @("7.252",     "  -0.43%  ")
    28 = @("16.12",     "  +0.39%  ")
    29 = @("1.408",     "  -0.48%  ")
    30 = @("0.05966",   "  -1.66%  ")
    31 = @("1.272",     "  -0.80%  ")
    32 = @("3.502",     "  -1.07%  ")
    33 = @("3.234",     "  -3.75%  ")
    34 = @("1.549",     "  -6.29%  ")
    35 = @("0.9465",    "  -3.85%  ")
    36 = @("2.412",     "  +0.03%  ")
    37 = @("2.754",     "  -0.60%  ")
    38 = @("0.5692",    "  -4.23%  ")
    39 = @("0.01598",   "  +0.08%  ")
    40 = @($null,       "  -1.25%  ")
    41 = @("0.8452",    "  -2.03%  ")
    42 = @($null,       "  -0.11%  ")
    43 = @("100.87",    "  +1.12%  ")
    44 = @("1.005.55",  "  -4.97%  ")
    45 = @("1.795.05",  "  +0.01%  ")
    46 = @("56.69",     "  -1.22%  ")
    47 = @($null,       "  -1.26%  ")
    48 = @("0.9986",    "  -0.53%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) {
        Set-TextValue $ws "D$row" $vals[0]
    }
    if ($null -ne $vals[1]) {
        Set-TextValue $ws "E$row" $vals[1]
    }
}

# Rows 49/50 swapped rank order (EnergySwap <-> Mantle) with refreshed
# price/volume figures.
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D49" "0.4304"
Set-TextValue $ws "E49" "  +1.71%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D50" "7.932"
Set-TextValue $ws "E50" "  -1.74%  "

# Row 51 price/volume refresh.
Set-TextValue $ws "D51" "1.481"
Set-TextValue $ws "E51" "  +1.10%  "
